$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 2 so existing data (rows 2-8) shifts down to rows 4-10
$ws.Range("A2:B3").Insert(-4121, 1)  # xlShiftDown, xlFormatFromRightOrBelow

# Fill in the new row 2 (date 42004, value 300792000000)
$ws.Cells.Item(2, 1).Value = 42004
$ws.Cells.Item(2, 2).Value = 300792000000

# Fill in the new row 3 (date 42369, value 281374000000)
$ws.Cells.Item(3, 1).Value = 42369
$ws.Cells.Item(3, 2).Value = 281374000000

# Copy style from an existing date cell (A4, which was originally A2) to the new A2:A3 cells
$ws.Range("A4").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# The insert operation may have carried formatting into column B; clear it to match plain cells
$ws.Range("B2:B3").ClearFormats()
